$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 20000
$ws.Range("J3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("N3").Value = -20228

$ws.Range("H28").Value = 78731.08
$ws.Range("I28").Value = 168302.17
$ws.Range("K28").Value = 168302.17
$ws.Range("M28").Value = -167817.17

$ws.Range("H33").Value = 235.45454
$ws.Range("I33").Value = 235.45454
$ws.Range("K33").Value = 235.45454
$ws.Range("M33").Value = -6.454540000000009

$ws.Range("H58").Value = 7168.6
$ws.Range("J58").Value = 9630.909
$ws.Range("L58").Value = 28892.727
$ws.Range("N58").Value = -29192.727

$ws.Range("H70").Value = 92861.55
$ws.Range("J70").Value = 126887.125
$ws.Range("L70").Value = 380661.375
$ws.Range("N70").Value = -381201.375

$ws.Range("H73").Value = 92861.55
$ws.Range("J73").Value = 126887.125
$ws.Range("L73").Value = 380661.375
$ws.Range("N73").Value = -382533.375

$ws.Range("H82").Value = 1767
$ws.Range("I82").Value = 1767
$ws.Range("K82").Value = 5301
$ws.Range("M82").Value = -4895

$ws.Range("H85").Value = 1767
$ws.Range("I85").Value = 1767
$ws.Range("K85").Value = 5301
$ws.Range("M85").Value = -3897

$ws.Range("H99").Value = 1146.2
$ws.Range("I99").Value = 184.66667
$ws.Range("J99").Value = 2588.5
$ws.Range("K99").Value = 554.00001
$ws.Range("L99").Value = 7765.5
$ws.Range("M99").Value = 943.99999
$ws.Range("N99").Value = -10761.5

$ws.Range("H100").Value = 6432.9414
$ws.Range("I100").Value = 2258.3845
$ws.Range("K100").Value = 2258.3845
$ws.Range("M100").Value = -1717.3845

$ws.Range("H101").Value = 704.4286
$ws.Range("J101").Value = 998
$ws.Range("L101").Value = 2994
$ws.Range("N101").Value = -6238

$ws.Range("H102").Value = 20000
$ws.Range("J102").Value = 20000
$ws.Range("L102").Value = 20000
$ws.Range("N102").Value = -26490

$ws.Range("H104").Value = 1066.3334
$ws.Range("I104").Value = 1066.3334
$ws.Range("K104").Value = 3199.0002
$ws.Range("M104").Value = -1452.0002

$ws.Range("H118").Value = 48124.75
$ws.Range("J118").Value = 999.5
$ws.Range("L118").Value = 2998.5
$ws.Range("N118").Value = -6312.5

$ws.Range("H127").Value = 11172.5
$ws.Range("I127").Value = 13697
$ws.Range("K127").Value = 41091
$ws.Range("M127").Value = -36131

$ws.Range("H129").Value = 15000.529
$ws.Range("I129").Value = 11608.909
$ws.Range("K129").Value = 34826.727
$ws.Range("M129").Value = -29826.727

$ws.Range("H132").Value = 2196.6924
$ws.Range("I132").Value = 1368.8182
$ws.Range("K132").Value = 4106.4546
$ws.Range("M132").Value = -1576.4546

$ws.Range("H138").Value = 4215.317
$ws.Range("I138").Value = 2477.423
$ws.Range("J138").Value = 5544.294
$ws.Range("K138").Value = 7432.268999999999
$ws.Range("L138").Value = 16632.882
$ws.Range("M138").Value = -2292.268999999999
$ws.Range("N138").Value = -26912.882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3690.1638
$ws.Range("I32").Value = 3072.9456
$ws.Range("K32").Value = 3072.9456
$ws.Range("M32").Value = -2785.9456

$ws.Range("H101").Value = 68000
$ws.Range("J101").Value = 68000
$ws.Range("L101").Value = 68000
$ws.Range("N101").Value = -74490

$ws.Range("H132").Value = 6411.7905
$ws.Range("I132").Value = 3764.742
$ws.Range("K132").Value = 11294.226
$ws.Range("M132").Value = -8764.226000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4436
$ws.Range("I86").Value = 4378.9287
$ws.Range("J86").Value = 4702.3335
$ws.Range("K86").Value = 4378.9287
$ws.Range("L86").Value = 4702.3335
$ws.Range("M86").Value = -3255.9287
$ws.Range("N86").Value = -6948.3335

$ws.Range("H89").Value = 4436
$ws.Range("I89").Value = 4378.9287
$ws.Range("J89").Value = 4702.3335
$ws.Range("K89").Value = 21894.6435
$ws.Range("L89").Value = 23511.6675
$ws.Range("M89").Value = -16278.6435
$ws.Range("N89").Value = -34743.6675

$ws.Range("H99").Value = 4527.273
$ws.Range("I99").Value = 4533.4443
$ws.Range("K99").Value = 4533.4443
$ws.Range("M99").Value = -3035.4443

$ws.Range("H134").Value = 5400.1113
$ws.Range("J134").Value = 7999.857
$ws.Range("L134").Value = 23999.571
$ws.Range("N134").Value = -29069.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4922.5557
$ws.Range("I62").Value = 4349.25
$ws.Range("J62").Value = 5381.2
$ws.Range("K62").Value = 4349.25
$ws.Range("L62").Value = 5381.2
$ws.Range("M62").Value = -3725.25
$ws.Range("N62").Value = -6629.2

$ws.Range("H65").Value = 4922.5557
$ws.Range("I65").Value = 4349.25
$ws.Range("J65").Value = 5381.2
$ws.Range("K65").Value = 21746.25
$ws.Range("L65").Value = 26906
$ws.Range("M65").Value = -18626.25
$ws.Range("N65").Value = -33146

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 41666764
$ws.Range("I11").Value = 41666764
$ws.Range("K11").Value = 125000292
$ws.Range("M11").Value = -125000152

$ws.Range("H33").Value = 2314839.5
$ws.Range("I33").Value = 2469155.5
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 14814933
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -14814650
$ws.Range("N33").Value = -1166

$ws.Range("H113").Value = 12345845
$ws.Range("I113").Value = 12345845
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 37037535
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -37035365
$ws.Range("N113").ClearContents()

$ws.Range("H121").Value = 16350459
$ws.Range("I121").Value = 727.7778
$ws.Range("K121").Value = 2183.3334
$ws.Range("M121").Value = -873.3334

$ws.Range("H131").Value = 5155.5713
$ws.Range("I131").Value = 1146.875
$ws.Range("J131").Value = 7622.4614
$ws.Range("K131").Value = 3440.625
$ws.Range("L131").Value = 22867.3842
$ws.Range("M131").Value = 1599.375
$ws.Range("N131").Value = -32947.3842

$ws.Range("H137").Value = 3002.75
$ws.Range("I137").Value = 2337
$ws.Range("K137").Value = 7011
$ws.Range("M137").Value = -1911

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 772281.5600000001
$ws.Range("I80").Value = 591265.4399999999
$ws.Range("K80").Value = 591265.4399999999
$ws.Range("M80").Value = -590267.4399999999

$ws.Range("H83").Value = 772281.5600000001
$ws.Range("I83").Value = 591265.4399999999
$ws.Range("K83").Value = 2956327.2
$ws.Range("M83").Value = -2951335.2

$ws.Range("H102").Value = 8577.579
$ws.Range("I102").Value = 8706.5
$ws.Range("K102").Value = 8706.5
$ws.Range("M102").Value = -7084.5

$ws.Range("H113").Value = 732015
$ws.Range("I113").Value = 1272302.4
$ws.Range("J113").Value = 11631.833
$ws.Range("K113").Value = 1272302.4
$ws.Range("L113").Value = 11631.833
$ws.Range("M113").Value = -1270132.4
$ws.Range("N113").Value = -15971.833

$ws.Range("H132").Value = 1672015
$ws.Range("J132").Value = 6800
$ws.Range("L132").Value = 20400
$ws.Range("N132").Value = -25460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6005250.5
$ws.Range("J40").Value = 7700
$ws.Range("L40").Value = 7700
$ws.Range("N40").Value = -7972

$ws.Range("H100").Value = 11737.375
$ws.Range("J100").Value = 25333
$ws.Range("L100").Value = 25333
$ws.Range("N100").Value = -26415

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12046.363
$ws.Range("I81").Value = 1081.8572
$ws.Range("K81").Value = 2163.7144
$ws.Range("M81").Value = -1102.7144

$ws.Range("H84").Value = 12046.363
$ws.Range("I84").Value = 1081.8572
$ws.Range("K84").Value = 10818.572
$ws.Range("M84").Value = -5514.572
